$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "59.164.54"
$ws.Cells.Item(2,5).Value = "  -3.07%  "
$ws.Cells.Item(3,4).Value = "2.653.74"
$ws.Cells.Item(3,5).Value = "  -1.16%  "
$ws.Cells.Item(4,5).Value = "  +0.18%  "
$ws.Cells.Item(5,4).Value = "'524.48"
$ws.Cells.Item(5,5).Value = "  +0.53%  "
$ws.Cells.Item(6,4).Value = "'144.76"
$ws.Cells.Item(6,5).Value = "  -1.92%  "
$ws.Cells.Item(7,5).Value = "  +0.28%  "
$ws.Cells.Item(8,4).Value = "'0.571"
$ws.Cells.Item(9,4).Value = "'7.03"
$ws.Cells.Item(9,5).Value = "  +9.65%  "
$ws.Cells.Item(10,5).Value = "  -3.51%  "
$ws.Cells.Item(11,4).Value = "'0.334"
$ws.Cells.Item(11,5).Value = "  -2.28%  "
$ws.Cells.Item(12,4).Value = "'0.131"
$ws.Cells.Item(12,5).Value = "  +1.61%  "
$ws.Cells.Item(13,4).Value = "3.119.22"
$ws.Cells.Item(13,5).Value = "  -1.32%  "
$ws.Cells.Item(14,4).Value = "59.203.27"
$ws.Cells.Item(14,5).Value = "  -2.91%  "
$ws.Cells.Item(15,4).Value = "'21.13"
$ws.Cells.Item(15,5).Value = "  -1.50%  "
$ws.Cells.Item(16,5).Value = "  -2.13%  "
$ws.Cells.Item(17,4).Value = "2.649.35"
$ws.Cells.Item(17,5).Value = "  -4.67%  "
$ws.Cells.Item(18,4).Value = "'340.70"
$ws.Cells.Item(18,5).Value = "  -3.88%  "
$ws.Cells.Item(19,5).Value = "  -4.29%  "
$ws.Cells.Item(20,4).Value = "'10.39"
$ws.Cells.Item(20,5).Value = "  -1.61%  "
$ws.Cells.Item(21,5).Value = "  +0.27%  "
$ws.Cells.Item(22,4).Value = "'0.999"
$ws.Cells.Item(22,5).Value = "  -0.02%  "
$ws.Cells.Item(23,4).Value = "'64.43"
$ws.Cells.Item(23,5).Value = "  +2.14%  "
$ws.Cells.Item(24,4).Value = "'0.420"
$ws.Cells.Item(24,5).Value = "  -0.94%  "
$ws.Cells.Item(25,5).Value = "  -2.07%  "
$ws.Cells.Item(26,4).Value = "'0.999"
$ws.Cells.Item(26,5).Value = "  +0.53%  "
$ws.Cells.Item(27,4).Value = "0.0₃0804"
$ws.Cells.Item(27,5).Value = "  -2.35%  "
$ws.Cells.Item(28,4).Value = "'7.12"
$ws.Cells.Item(28,5).Value = "  -2.74%  "
$ws.Cells.Item(29,4).Value = "'6.68"
$ws.Cells.Item(29,5).Value = "  -2.51%  "
$ws.Cells.Item(30,4).Value = "'0.998"
$ws.Cells.Item(30,5).Value = "  +0.04%  "
$ws.Cells.Item(31,5).Value = "  -0.18%  "
$ws.Cells.Item(32,4).Value = "'18.86"
$ws.Cells.Item(32,5).Value = "  -1.64%  "
$ws.Cells.Item(33,4).Value = "'149.27"
$ws.Cells.Item(33,5).Value = "  -0.19%  "
$ws.Cells.Item(34,5).Value = "  -2.56%  "
$ws.Cells.Item(35,5).Value = "  -3.27%  "
$ws.Cells.Item(36,4).Value = "'0.898"
$ws.Cells.Item(36,5).Value = "  -5.51%  "
$ws.Cells.Item(37,4).Value = "'0.883"
$ws.Cells.Item(37,5).Value = "  +0.92%  "
$ws.Cells.Item(38,4).Value = "'36.73"
$ws.Cells.Item(38,5).Value = "  +0.05%  "
$ws.Cells.Item(39,5).Value = "  -5.76%  "
$ws.Cells.Item(40,5).Value = "  -3.67%  "
$ws.Cells.Item(41,4).Value = "'0.617"
$ws.Cells.Item(41,5).Value = "  +0.59%  "
$ws.Cells.Item(42,2).Value = "EnergySwap"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(42,4).Value = "'20.09"
$ws.Cells.Item(42,5).Value = "  +0.17%  "
$ws.Cells.Item(43,2).Value = "FirstDigitalUSD"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(43,4).Value = "'0.999"
$ws.Cells.Item(43,5).Value = "  +0.10%  "
$ws.Cells.Item(44,2).Value = "Bittensor"
$ws.Cells.Item(44,3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(44,4).Value = "'275.16"
$ws.Cells.Item(44,5).Value = "  -3.31%  "
$ws.Cells.Item(45,4).Value = "'0.0971"
$ws.Cells.Item(45,5).Value = "  -2.24%  "
$ws.Cells.Item(46,4).Value = "'0.0533"
$ws.Cells.Item(46,5).Value = "  -1.56%  "
$ws.Cells.Item(47,5).Value = "  +1.79%  "
$ws.Cells.Item(48,4).Value = "2.031.88"
$ws.Cells.Item(48,5).Value = "  -5.02%  "
$ws.Cells.Item(49,4).Value = "'4.78"
$ws.Cells.Item(49,5).Value = "  -2.17%  "
$ws.Cells.Item(50,5).Value = "  -2.85%  "
$ws.Cells.Item(51,4).Value = "'18.87"
